$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# hunk 0
$ws.Range("H9").Value = 600.875
$ws.Range("J9").Value = 1359
$ws.Range("L9").Value = 1359
$ws.Range("N9").Value = -1697
# hunk 1
$ws.Range("H33").Value = 10000618
$ws.Range("I33").Value = 15625473
$ws.Range("J33").Value = 875.44446
$ws.Range("K33").Value = 15625473
$ws.Range("L33").Value = 875.44446
$ws.Range("M33").Value = -15625244
$ws.Range("N33").Value = -1333.44446
# hunk 2
$ws.Range("H74").Value = 9210.333000000001
$ws.Range("I74").Value = 9190.5
$ws.Range("K74").Value = 9190.5
$ws.Range("M74").Value = -8254.5
# hunk 3
$ws.Range("H77").Value = 9210.333000000001
$ws.Range("I77").Value = 9190.5
$ws.Range("K77").Value = 45952.5
$ws.Range("M77").Value = -41272.5
# hunk 4
$ws.Range("H111").Value = 1085.375
$ws.Range("I111").Value = 954.7143
$ws.Range("K111").Value = 2864.1429
$ws.Range("M111").Value = 202.8571000000002
# hunk 5
$ws.Range("H113").Value = 8320.315000000001
$ws.Range("I113").Value = 11379.8
$ws.Range("J113").Value = 4920.8887
$ws.Range("K113").Value = 11379.8
$ws.Range("L113").Value = 4920.8887
$ws.Range("M113").Value = -8125.799999999999
$ws.Range("N113").Value = -11428.8887
# hunk 6
$ws.Range("H129").Value = 1345.6957
$ws.Range("I129").Value = 713.1667
$ws.Range("J129").Value = 1568.9412
$ws.Range("K129").Value = 2139.5001
$ws.Range("L129").Value = 4706.8236
$ws.Range("M129").Value = 2860.4999
$ws.Range("N129").Value = -14706.8236
# hunk 7
$ws.Range("H138").Value = 4072.054
$ws.Range("I138").Value = 5329.4443
$ws.Range("K138").Value = 15988.3329
$ws.Range("M138").Value = -10848.3329

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# hunk 8
$ws.Range("H2").Value = 3726.5
$ws.Range("I2").Value = 3929.3076
$ws.Range("J2").Value = 3199.2
$ws.Range("K2").Value = 3929.3076
$ws.Range("L2").Value = 3199.2
$ws.Range("M2").Value = -3816.3076
$ws.Range("N2").Value = -3425.2
# hunk 9
$ws.Range("H61").Value = 430151.7
$ws.Range("I61").Value = 2382.0833
$ws.Range("J61").Value = 3281949
$ws.Range("K61").Value = 2382.0833
$ws.Range("L61").Value = 3281949
$ws.Range("M61").Value = -2170.0833
$ws.Range("N61").Value = -3282373
# hunk 10
$ws.Range("H74").Value = 6325.553
$ws.Range("I74").Value = 2103.7856
$ws.Range("J74").Value = 41788.4
$ws.Range("K74").Value = 2103.7856
$ws.Range("L74").Value = 41788.4
$ws.Range("M74").Value = -1229.7856
$ws.Range("N74").Value = -43536.4
# hunk 11
$ws.Range("H77").Value = 6325.553
$ws.Range("I77").Value = 2103.7856
$ws.Range("J77").Value = 41788.4
$ws.Range("K77").Value = 10518.928
$ws.Range("L77").Value = 208942
$ws.Range("M77").Value = -6150.928
$ws.Range("N77").Value = -217678
# hunk 12
$ws.Range("H102").Value = 5251.5
$ws.Range("I102").Value = 5709
$ws.Range("J102").Value = 2964
$ws.Range("K102").Value = 5709
$ws.Range("L102").Value = 2964
$ws.Range("M102").Value = -4087
$ws.Range("N102").Value = -6208
# hunk 13
$ws.Range("H116").Value = 3726.5
$ws.Range("I116").Value = 3929.3076
$ws.Range("J116").Value = 3199.2
$ws.Range("K116").Value = 3929.3076
$ws.Range("L116").Value = 3199.2
$ws.Range("M116").Value = -1635.3076
$ws.Range("N116").Value = -7787.2
# hunk 14
$ws.Range("H132").Value = 450651.75
$ws.Range("I132").Value = 2528.224
$ws.Range("J132").Value = 2307163.5
$ws.Range("K132").Value = 7584.672
$ws.Range("L132").Value = 6921490.5
$ws.Range("M132").Value = -5054.672
$ws.Range("N132").Value = -6926550.5
# hunk 15
$ws.Range("H136").Value = 430151.7
$ws.Range("I136").Value = 2382.0833
$ws.Range("J136").Value = 3281949
$ws.Range("K136").Value = 7146.249899999999
$ws.Range("L136").Value = 9845847
$ws.Range("M136").Value = -4596.249899999999
$ws.Range("N136").Value = -9850947

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# hunk 16
$ws.Range("H3").Value = 3726.5
$ws.Range("I3").Value = 3929.3076
$ws.Range("J3").Value = 3199.2
$ws.Range("K3").Value = 3929.3076
$ws.Range("L3").Value = 3199.2
$ws.Range("M3").Value = -3815.3076
$ws.Range("N3").Value = -3427.2
# hunk 17
$ws.Range("H99").Value = 2767.3333
$ws.Range("I99").Value = 2197.5
$ws.Range("J99").Value = 7326
$ws.Range("K99").Value = 2197.5
$ws.Range("L99").Value = 7326
$ws.Range("M99").Value = -699.5
$ws.Range("N99").Value = -10322
# hunk 18
$ws.Range("H105").Value = 4677.625
$ws.Range("I105").Value = 4482
$ws.Range("J105").Value = 5003.6665
$ws.Range("K105").Value = 4482
$ws.Range("L105").Value = 5003.6665
$ws.Range("M105").Value = -2735
$ws.Range("N105").Value = -8497.666499999999
# hunk 19
$ws.Range("H134").Value = 11053
$ws.Range("I134").Value = 7886.057
$ws.Range("K134").Value = 23658.171
$ws.Range("M134").Value = -21123.171

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# hunk 20
$ws.Range("H3").Value = 616.6667
$ws.Range("I3").Value = 616.6667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 616.6667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -503.6667
$ws.Range("N3").ClearContents()
# hunk 21
$ws.Range("H99").Value = 3715.4546
$ws.Range("I99").Value = 3696.375
$ws.Range("J99").Value = 3766.3333
$ws.Range("K99").Value = 3696.375
$ws.Range("L99").Value = 3766.3333
$ws.Range("M99").Value = -2198.375
$ws.Range("N99").Value = -6762.3333
# hunk 22
$ws.Range("H126").Value = 3715.4546
$ws.Range("I126").Value = 3696.375
$ws.Range("J126").Value = 3766.3333
$ws.Range("K126").Value = 11089.125
$ws.Range("L126").Value = 11298.9999
$ws.Range("M126").Value = -8619.125
$ws.Range("N126").Value = -16238.9999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# hunk 23
$ws.Range("H7").Value = 178.6
$ws.Range("J7").Value = 303.33334
$ws.Range("L7").Value = 910.0000200000001
$ws.Range("N7").Value = -1134.00002
# hunk 24
$ws.Range("H12").Value = 29.5
$ws.Range("I12").Value = 26.75
$ws.Range("J12").Value = 32.25
$ws.Range("K12").Value = 80.25
$ws.Range("L12").Value = 96.75
$ws.Range("M12").Value = 92.75
$ws.Range("N12").Value = -442.75
# hunk 25
$ws.Range("H34").Value = 2190.125
$ws.Range("I34").Value = 1113.0714
$ws.Range("J34").Value = 3698
$ws.Range("K34").Value = 3339.2142
$ws.Range("L34").Value = 11094
$ws.Range("M34").Value = -3255.2142
$ws.Range("N34").Value = -11262
# hunk 26
$ws.Range("H39").Value = 4547.5293
$ws.Range("I39").Value = 2162.4
$ws.Range("J39").Value = 5541.3335
$ws.Range("K39").Value = 6487.200000000001
$ws.Range("L39").Value = 16624.0005
$ws.Range("M39").Value = -6193.200000000001
$ws.Range("N39").Value = -17212.0005
# hunk 27
$ws.Range("H55").Value = 3446.3333
$ws.Range("I55").Value = 3446.3333
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 10338.9999
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -10161.9999
$ws.Range("N55").ClearContents()
# hunk 28
$ws.Range("H129").Value = 19390.455
$ws.Range("I129").Value = 50682.5
$ws.Range("J129").Value = 1509.2858
$ws.Range("K129").Value = 152047.5
$ws.Range("L129").Value = 4527.857400000001
$ws.Range("M129").Value = -147047.5
$ws.Range("N129").Value = -14527.8574
# hunk 29
$ws.Range("H131").Value = 5784.4507
$ws.Range("J131").Value = 5239.5977
$ws.Range("L131").Value = 15718.7931
$ws.Range("N131").Value = -25798.7931

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# hunk 30
$ws.Range("H40").Value = 4382.5
$ws.Range("I40").Value = 3332
$ws.Range("J40").Value = 5433
$ws.Range("K40").Value = 3332
$ws.Range("L40").Value = 5433
$ws.Range("M40").Value = -3196
$ws.Range("N40").Value = -5705
# hunk 31
$ws.Range("H55").Value = 3099.2222
$ws.Range("I55").Value = 3135
$ws.Range("K55").Value = 3135
$ws.Range("M55").Value = -2962
# hunk 32
$ws.Range("H122").Value = 10329
$ws.Range("I122").Value = 11333
$ws.Range("K122").Value = 33999
$ws.Range("M122").Value = -31549

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# hunk 33
$ws.Range("H136").Value = 232818
$ws.Range("I136").Value = 1562.7916
$ws.Range("K136").Value = 4688.3748
$ws.Range("M136").Value = -2138.3748

